$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '24.640.95'
$ws.Range('E2').Value = '  -0.08%  '

# Row 3
$ws.Range('D3').Value = '1.700.96'
$ws.Range('E3').Value = '  +0.32%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9970'
$ws.Range('E4').Value = '  -1.02%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '314.22'
$ws.Range('E5').Value = '  -1.44%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9947'
$ws.Range('E6').Value = '  -1.10%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3961'
$ws.Range('E7').Value = '  -0.34%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4063'
$ws.Range('E8').Value = '  +0.93%  '

# Row 9
$ws.Range('B9').Value = 'BinanceUSD'
$ws.Range('C9').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.9933'
$ws.Range('E9').Value = '  -1.05%  '

# Row 10
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.511'
$ws.Range('E10').Value = '  +5.40%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.14'
$ws.Range('E11').Value = '  +8.57%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08778'
$ws.Range('E12').Value = '  -0.62%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.316'
$ws.Range('E13').Value = '  +8.48%  '

# Row 14
$ws.Range('E14').Value = '  -0.63%  '

# Row 15
$ws.Range('E15').Value = '  -0.77%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.480'
$ws.Range('E16').Value = '  +2.85%  '

# Row 17
$ws.Range('D17').Value = '1.701.80'
$ws.Range('E17').Value = '  +0.40%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '100.82'
$ws.Range('E18').Value = '  -1.81%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07056'
$ws.Range('E19').Value = '  +2.62%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.46'
$ws.Range('E20').Value = '  -1.69%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.753'
$ws.Range('E21').Value = '  -1.49%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9941'
$ws.Range('E22').Value = '  -1.19%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.20'
$ws.Range('E23').Value = '  +0.84%  '

# Row 24
$ws.Range('D24').Value = '24.703.35'
$ws.Range('E24').Value = '  +0.21%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.979'
$ws.Range('E25').Value = '  +3.21%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.306'
$ws.Range('E26').Value = '  -0.41%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.40'
$ws.Range('E27').Value = '  +0.22%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '158.42'
$ws.Range('E28').Value = '  -0.96%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.125'
$ws.Range('E29').Value = '  -3.44%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '133.04'
$ws.Range('E30').Value = '  -0.58%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.434'
$ws.Range('E31').Value = '  +24.72%  '

# Row 32
$ws.Range('D32').Value = '1.886.10'
$ws.Range('E32').Value = '  +0.21%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.089'
$ws.Range('E33').Value = '  -9.41%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.436'
$ws.Range('E34').Value = '  +20.59%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.08687'
$ws.Range('E35').Value = '  -3.52%  '

# Row 36
$ws.Range('E36').Value = '  +0.59%  '

# Row 37
$ws.Range('B37').Value = 'WEMIXTOKEN'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.941'
$ws.Range('E37').Value = '  +3.23%  '

# Row 38
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2731'
$ws.Range('E38').Value = '  -0.25%  '

# Row 39
$ws.Range('E39').Value = '  -4.73%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.02763'
$ws.Range('E40').Value = '  +7.60%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.08976'
$ws.Range('E41').Value = '  +0.42%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.481'
$ws.Range('E42').Value = '  +0.43%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.7653'
$ws.Range('E43').Value = '  -0.38%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.7225'
$ws.Range('E44').Value = '  +0.04%  '

# Row 45
$ws.Range('E45').Value = '  +0.62%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.461'
$ws.Range('E46').Value = '  -1.23%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.160'
$ws.Range('E47').Value = '  -0.20%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.9932'
$ws.Range('E48').Value = '  -1.18%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '141.34'
$ws.Range('E49').Value = '  -1.35%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.317'
$ws.Range('E50').Value = '  +12.83%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.08020'
